$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (TruckID, AssignedDockPosition, start_loading_time, end_loading_time)
# for rows 2-12, replacing the old rows 2-10.
$data = @(
    @(4, 1, 5, 5),
    @(1, 2, 5, 5),
    @(5, 2, 10, 10),
    @(6, 2, 15, 15),
    @(7, 2, 20, 20),
    @(9, 2, 25, 25),
    @(10, 2, 30, 31),
    @(2, 3, 5, 5),
    @(3, 3, 10, 10),
    @(8, 3, 15, 15),
    @(9, 3, 20, 20)
)

$row = 2
foreach ($rowData in $data) {
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $row++
}
